# Updated to larger Llama model
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows of evaluation results
$ws.Range("A5").Value = "BM25 Retriever + Semantic chunking"
$ws.Range("B5").Value = 0.7338
$ws.Range("C5").Value = 0.533
$ws.Range("D5").Value = 0.1821

$ws.Range("A6").Value = "BM25 Retriever + Semantic chunking + Llama 3.2:3B"
$ws.Range("B6").Value = 0.7675
$ws.Range("C6").Value = 0.5722
$ws.Range("D6").Value = 0.2478

$ws.Range("A7").Value = "BM25 Retriever + Semantic chunking + Llama 3.2:3B + Prompt Templates"

# Widen column A to fit the longer labels (engine adds a fixed padding
# offset of 5/6 character when storing ColumnWidth, so back it out here
# to land on a stored width of exactly 65)
$ws.Columns.Item(1).ColumnWidth = 64.1666667

# Move the active selection to the next empty row, like Excel does after data entry
$ws.Range("A8").Select()
